# Update TPM-derived NATMI metrics for the Wnt6-Fzd7 ligand-receptor pair sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending=FAPs, Target=ECs
$ws.Range("I2").Value = 0.2893336272138922
$ws.Range("J2").Value = 0.2893336272138922
$ws.Range("M2").Value = 0.6068319999999999
$ws.Range("N2").Value = 1.820496
$ws.Range("O2").Value = 0.03392274820144286
$ws.Range("P2").Value = 0.03392274820144286
$ws.Range("Q2").Value = 0.2825241901813333
$ws.Range("R2").Value = 2.542717711632
$ws.Range("S2").Value = 0.009814991782187003
$ws.Range("T2").Value = 0.009814991782187001

# Row 3: Sending=FAPs, Target=FAPs
$ws.Range("I3").Value = 0.2893336272138922
$ws.Range("J3").Value = 0.2893336272138922
$ws.Range("O3").Value = 0.4504903529585388
$ws.Range("P3").Value = 0.4504903529585388
$ws.Range("S3").Value = 0.1303420078463606
$ws.Range("T3").Value = 0.1303420078463606

# Row 4: Sending=FAPs, Target=MuSCs
$ws.Range("I4").Value = 0.2893336272138922
$ws.Range("J4").Value = 0.2893336272138922
$ws.Range("M4").Value = 9.223151
$ws.Range("N4").Value = 27.669453
$ws.Range("O4").Value = 0.5155868988400183
$ws.Range("P4").Value = 0.5155868988400183
$ws.Range("Q4").Value = 4.294043931755668
$ws.Range("R4").Value = 38.64639538580101
$ws.Range("S4").Value = 0.1491766275853446
$ws.Range("T4").Value = 0.1491766275853446

# Row 5: Sending=MuSCs, Target=ECs
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.143547
$ws.Range("H5").Value = 3.430641
$ws.Range("I5").Value = 0.7106663727861078
$ws.Range("J5").Value = 0.7106663727861078
$ws.Range("M5").Value = 0.6068319999999999
$ws.Range("N5").Value = 1.820496
$ws.Range("O5").Value = 0.03392274820144286
$ws.Range("P5").Value = 0.03392274820144286
$ws.Range("Q5").Value = 0.6939409131039999
$ws.Range("R5").Value = 6.245468217935999
$ws.Range("S5").Value = 0.02410775641925586
$ws.Range("T5").Value = 0.02410775641925586

# Row 6: Sending=MuSCs, Target=FAPs
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.143547
$ws.Range("H6").Value = 3.430641
$ws.Range("I6").Value = 0.7106663727861078
$ws.Range("J6").Value = 0.7106663727861078
$ws.Range("O6").Value = 0.4504903529585388
$ws.Range("P6").Value = 0.4504903529585388
$ws.Range("Q6").Value = 9.215458754114
$ws.Range("R6").Value = 82.93912878702599
$ws.Range("S6").Value = 0.3201483451121782
$ws.Range("T6").Value = 0.3201483451121782

# Row 7: Sending=MuSCs, Target=MuSCs
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.143547
$ws.Range("H7").Value = 3.430641
$ws.Range("I7").Value = 0.7106663727861078
$ws.Range("J7").Value = 0.7106663727861078
$ws.Range("M7").Value = 9.223151
$ws.Range("N7").Value = 27.669453
$ws.Range("O7").Value = 0.5155868988400183
$ws.Range("P7").Value = 0.5155868988400183
$ws.Range("Q7").Value = 10.547106656597
$ws.Range("R7").Value = 94.923959909373
$ws.Range("S7").Value = 0.3664102712546737
$ws.Range("T7").Value = 0.3664102712546737
